$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.944.03'
$ws.Range('E2').Value = '  +0.34%  '
$ws.Range('D3').Value = '1.636.97'
$ws.Range('E3').Value = '  -0.21%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = '212.21'
$ws.Range('E5').Value = '  -0.28%  '
$ws.Range('E6').Value = '  -0.25%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').Value = '23.28'
$ws.Range('E8').Value = '  -0.89%  '
$ws.Range('E9').Value = '  -2.45%  '
$ws.Range('E10').Value = '  +0.03%  '
$ws.Range('D11').Value = '0.0881'
$ws.Range('E11').Value = '  +1.46%  '
$ws.Range('D12').Value = '1.870.82'
$ws.Range('E12').Value = '  -0.11%  '
$ws.Range('D13').Value = '1.645.95'
$ws.Range('E13').Value = '  +0.30%  '
$ws.Range('E14').Value = '  +0.05%  '
$ws.Range('D15').Value = '0.567'
$ws.Range('E15').Value = '  +0.76%  '
$ws.Range('D16').Value = '65.37'
$ws.Range('E16').Value = '  -0.42%  '
$ws.Range('D17').Value = '27.951.42'
$ws.Range('E17').Value = '  +0.33%  '
$ws.Range('D18').Value = '231.24'
$ws.Range('E18').Value = '  -0.16%  '
$ws.Range('E19').Value = '  -0.19%  '
$ws.Range('D20').Value = '7.53'
$ws.Range('E20').Value = '  -1.90%  '
$ws.Range('D22').Value = '10.39'
$ws.Range('E22').Value = '  -3.10%  '
$ws.Range('E23').Value = '  -0.39%  '
$ws.Range('E24').Value = '  -3.91%  '
$ws.Range('D25').Value = '153.27'
$ws.Range('E25').Value = '  +1.17%  '
$ws.Range('D26').Value = '6.95'
$ws.Range('E26').Value = '  +0.56%  '
$ws.Range('E27').Value = '  -0.27%  '
$ws.Range('D28').Value = '15.60'
$ws.Range('E28').Value = '  -0.53%  '
$ws.Range('E29').Value = '  +0.06%  '
$ws.Range('D30').Value = '1.18'
$ws.Range('E30').Value = '  -0.07%  '
$ws.Range('D31').Value = '0.0484'
$ws.Range('E31').Value = '  +0.18%  '
$ws.Range('D32').Value = '3.38'
$ws.Range('E32').Value = '  +1.98%  '
$ws.Range('D33').Value = '1.400.76'
$ws.Range('E33').Value = '  -3.94%  '
$ws.Range('D34').Value = '3.07'
$ws.Range('E34').Value = '  -1.28%  '
$ws.Range('E35').Value = '  +1.49%  '
$ws.Range('E36').Value = '  +1.42%  '
$ws.Range('E37').Value = '  +0.40%  '
$ws.Range('D38').Value = '0.562'
$ws.Range('E38').Value = '  -0.10%  '
$ws.Range('D39').Value = '0.927'
$ws.Range('E39').Value = '  +0.92%  '
$ws.Range('E40').Value = '  -1.39%  '
$ws.Range('E41').Value = '  +0.73%  '
$ws.Range('E42').Value = '  -0.03%  '
$ws.Range('D43').Value = '66.92'
$ws.Range('E43').Value = '  -3.30%  '
$ws.Range('D44').Value = '5.52'
$ws.Range('E44').Value = '  +2.61%  '
$ws.Range('E45').Value = '  +2.01%  '
$ws.Range('E46').Value = '  -0.58%  '
$ws.Range('D47').Value = '1.779.09'
$ws.Range('D48').Value = '87.93'
$ws.Range('E48').Value = '  -0.52%  '
$ws.Range('D49').Value = '0.0999'
$ws.Range('E49').Value = '  -0.77%  '
$ws.Range('E50').Value = '  -0.35%  '
$ws.Range('E51').Value = '  -1.65%  '
